$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 2 (old rows 2-4 shift down to 5-7)
$ws.Rows("2:4").Insert()
# The insert copies formatting from the row above; strip it back to Normal
# so the new rows match the unstyled data rows elsewhere in the sheet.
$ws.Rows("2:4").Style = "Normal"

# --- Row 2: Estudiantes respiran con el aumento del pasaje de TransMilenio (Alerta Bogotá) ---
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-13"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = "Estudiantes respiran con el aumento del pasaje de TransMilenio: no tendrán que gastarse lo del almuerzo"
$ws.Cells.Item(2, 3).Value = "Alerta Bogotá"
$ws.Cells.Item(2, 4).Value = "Bogotá"
$ws.Cells.Item(2, 5).Value = "https://www.alertabogota.com/noticias/local"
# F2 (resumen) is intentionally left blank for this article

# --- Row 3: Treinta niños quedaron sin aulas en zona rural de Nechí, Antioquia (Infobae) ---
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "2026-01-13"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Value = "Treinta niños quedaron sin aulas en zona rural de Nechí, Antioquia: incendio incineró la única escuela de un caserío"
$ws.Cells.Item(3, 3).Value = "Infobae"
$ws.Cells.Item(3, 4).Value = "Sin identificar"
$ws.Cells.Item(3, 5).Value = "https://www.infobae.com/colombia/2026/01/13/treinta-ninos-quedaron-sin-aulas-en-zona-rural-de-nechi-antioquia-incendio-incinero-la-unica-escuela-de-un-caserio/"
$ws.Cells.Item(3, 6).Value = "PorJuan Sánchez Romero"

# --- Row 4: Tres estudiantes resultan heridos tras caer desde un bus de dos niveles (Diario ADN) ---
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "2026-01-13"
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(4, 2).Value = "Tres estudiantes resultan heridos tras caer desde un bus de dos niveles; autoridades investigan"
$ws.Cells.Item(4, 3).Value = "Diario ADN"
$ws.Cells.Item(4, 4).Value = "Sin identificar"
$ws.Cells.Item(4, 5).Value = "https://www.diarioadn.co/seccion/actualidad"
$ws.Cells.Item(4, 6).Value = "Los menores, de 13 y 14 años, tuvieron que ser hospitalizados."
